$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.013.20"
$ws.Range("E2").Value = "  +5.18%  "
$ws.Range("D3").Value = "2.259.15"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.04"
$ws.Range("E5").Value = "  +3.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.30"
$ws.Range("E6").Value = "  +5.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.533"
$ws.Range("E7").Value = "  +3.71%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +3.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "54.72"
$ws.Range("E10").Value = "  +8.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.29"
$ws.Range("E11").Value = "  +5.81%  "
$ws.Range("E12").Value = "  +2.48%  "
$ws.Range("E13").Value = "  +2.79%  "
$ws.Range("E14").Value = "  +3.55%  "
$ws.Range("D15").Value = "2.608.76"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.10"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("D17").Value = "2.279.78"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("E18").Value = "  +3.42%  "
$ws.Range("D19").Value = "41.882.07"
$ws.Range("E19").Value = "  +5.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.18"
$ws.Range("E20").Value = "  +9.56%  "
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("E22").Value = "  +3.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.00"
$ws.Range("E23").Value = "  +1.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "241.16"
$ws.Range("E24").Value = "  +1.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.56"
$ws.Range("E25").Value = "  +4.06%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +3.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.88"
$ws.Range("E28").Value = "  +3.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.30"
$ws.Range("E29").Value = "  +13.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.65"
$ws.Range("E30").Value = "  +4.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.89"
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.88"
$ws.Range("E32").Value = "  +5.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.16"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0744"
$ws.Range("E35").Value = "  +4.20%  "
$ws.Range("E36").Value = "  +2.75%  "
$ws.Range("E37").Value = "  +2.76%  "
$ws.Range("E38").Value = "  +3.57%  "
$ws.Range("E39").Value = "  +4.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.52"
$ws.Range("E40").Value = "  +7.81%  "
$ws.Range("E41").Value = "  +2.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.92"
$ws.Range("E42").Value = "  +5.50%  "
$ws.Range("D43").Value = "2.055.41"
$ws.Range("E43").Value = "  -2.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.60"
$ws.Range("E44").Value = "  +7.95%  "
$ws.Range("E45").Value = "  +2.55%  "
$ws.Range("E46").Value = "  +2.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.88"
$ws.Range("E47").Value = "  +5.62%  "
$ws.Range("E48").Value = "  +2.67%  "
$ws.Range("E49").Value = "  +3.78%  "
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.58"
$ws.Range("E51").Value = "  +5.09%  "
